$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.003.04'
$ws.Range('E2').Value = '  -1.56%  '
$ws.Range('D3').Value = '1.822.35'
$ws.Range('E3').Value = '  -1.08%  '
$ws.Range('E4').Value = '  -0.65%  '
$ws.Range('D5').Value = '''309.72'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.73%  '
$ws.Range('D6').Value = '''1.008'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.55%  '
$ws.Range('D7').Value = '''0.4637'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.73%  '
$ws.Range('D8').Value = '''0.3643'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.81%  '
$ws.Range('D9').Value = '''0.07302'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.28%  '
$ws.Range('E10').Value = '  -2.31%  '
$ws.Range('D11').Value = '''19.85'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.20%  '
$ws.Range('D12').Value = '1.884.27'
$ws.Range('E12').Value = '  +1.57%  '
$ws.Range('D13').Value = '''0.07592'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.03%  '
$ws.Range('D14').Value = '''93.32'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.05%  '
$ws.Range('D15').Value = '''5.333'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.82%  '
$ws.Range('D16').Value = '''6.490'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.68%  '
$ws.Range('D17').Value = '''1.008'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.82%  '
$ws.Range('D18').Value = '''0.000008643'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.42%  '
$ws.Range('E19').Value = '  -0.62%  '
$ws.Range('D20').Value = '27.520.18'
$ws.Range('E20').Value = '  +0.30%  '
$ws.Range('D21').Value = '''14.48'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.40%  '
$ws.Range('E22').Value = '  -3.45%  '
$ws.Range('E23').Value = '  -1.55%  '
$ws.Range('D24').Value = '2.133.60'
$ws.Range('E24').Value = '  +2.68%  '
$ws.Range('D25').Value = '''151.93'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.40%  '
$ws.Range('D26').Value = '''1.859'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.66%  '
$ws.Range('D27').Value = '''18.25'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.20%  '
$ws.Range('D28').Value = '''2.094'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.49%  '
$ws.Range('B29').Value = 'BitcoinCash'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D29').Value = '''116.09'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.80%  '
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').Value = '''5.086'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.72%  '
$ws.Range('D31').Value = '''0.08903'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.94%  '
$ws.Range('D32').Value = '''2.952'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.01%  '
$ws.Range('D33').Value = '''0.7300'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.86%  '
$ws.Range('D34').Value = '''1.143'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.09%  '
$ws.Range('D35').Value = '''4.429'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.10%  '
$ws.Range('E36').Value = '  -0.58%  '
$ws.Range('D37').Value = '''2.508'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +5.50%  '
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').Value = '''0.05277'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.63%  '
$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D39').Value = '''1.075'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.97%  '
$ws.Range('D40').Value = '''0.01917'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.33%  '
$ws.Range('D41').Value = '''2.924'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.61%  '
$ws.Range('D42').Value = '''7.185'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.77%  '
$ws.Range('D43').Value = '''0.5223'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.54%  '
$ws.Range('D44').Value = '''0.1634'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.98%  '
$ws.Range('D45').Value = '''8.276'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.30%  '
$ws.Range('D46').Value = '''0.4867'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.27%  '
$ws.Range('E47').Value = '  -0.63%  '
$ws.Range('D48').Value = '''10.16'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.51%  '
$ws.Range('D49').Value = '''103.36'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').Value = '''1.635'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.02%  '
$ws.Range('D51').Value = '''0.06228'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.60%  '
